$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "shiaa_29"
$ws.Range("B2").Value = 29

$ws.Range("B2").Select()
